$wb = $excel.ActiveWorkbook

# 1. Add the new quarter worksheet "2022-Q3" right after "Zong Ji" (total) and before "2022-Q1"
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Add($null, $wsTotal)
$wsQ3.Name = "2022-Q3"

# Header row (copy the bold/boxed style used elsewhere in the workbook)
$wsQ3.Range("B1").Value = '基金代码'
$wsQ3.Range("C1").Value = '基金名称'
$wsQ3.Range("D1").Value = '基金规模'
$wsQ3.Range("E1").Value = '股票总仓位'
$wsQ3.Range("F1").Value = '仓位占比'
$wsQ3.Range("G1").Value = '持有市值(亿元)'
$wsQ3.Range("H1").Value = '仓位排名'

# Columns B-G hold text (fund codes / decimal strings) - format as Text first so
# leading/trailing zeros (e.g. "001511", "2.50", "0.0300") are preserved verbatim
$wsQ3.Range("B2:G21").NumberFormat = "@"

# Fund rows - column A (row index) and H (rank) are real numbers, B-G are text
$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = '001511'
$wsQ3.Range("C2").Value = '兴全新视野灵活配置定期开放混合'
$wsQ3.Range("D2").Value = '111.19'
$wsQ3.Range("E2").Value = '87.43'
$wsQ3.Range("F2").Value = '3.49'
$wsQ3.Range("G2").Value = '3.8805'
$wsQ3.Range("H2").Value = 6
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = '163415'
$wsQ3.Range("C3").Value = '兴全商业模式优选混合（LOF）'
$wsQ3.Range("D3").Value = '106.72'
$wsQ3.Range("E3").Value = '93.89'
$wsQ3.Range("F3").Value = '3.38'
$wsQ3.Range("G3").Value = '3.6071'
$wsQ3.Range("H3").Value = 6
$wsQ3.Range("A4").Value = 2
$wsQ3.Range("B4").Value = '011056'
$wsQ3.Range("C4").Value = '博时汇兴回报一年持有期灵活配置混合'
$wsQ3.Range("D4").Value = '96.44'
$wsQ3.Range("E4").Value = '52.92'
$wsQ3.Range("F4").Value = '2.50'
$wsQ3.Range("G4").Value = '2.4110'
$wsQ3.Range("H4").Value = 8
$wsQ3.Range("A5").Value = 3
$wsQ3.Range("B5").Value = '013797'
$wsQ3.Range("C5").Value = '博时优质鑫选一年持有期混合A'
$wsQ3.Range("D5").Value = '47.45'
$wsQ3.Range("E5").Value = '80.56'
$wsQ3.Range("F5").Value = '2.68'
$wsQ3.Range("G5").Value = '1.2717'
$wsQ3.Range("H5").Value = 10
$wsQ3.Range("A6").Value = 4
$wsQ3.Range("B6").Value = '001236'
$wsQ3.Range("C6").Value = '博时丝路主题股票A'
$wsQ3.Range("D6").Value = '13.40'
$wsQ3.Range("E6").Value = '83.01'
$wsQ3.Range("F6").Value = '2.62'
$wsQ3.Range("G6").Value = '0.3511'
$wsQ3.Range("H6").Value = 7
$wsQ3.Range("A7").Value = 5
$wsQ3.Range("B7").Value = '009740'
$wsQ3.Range("C7").Value = '博时研究臻选三年持有期灵活配置混合A'
$wsQ3.Range("D7").Value = '7.61'
$wsQ3.Range("E7").Value = '82.29'
$wsQ3.Range("F7").Value = '2.79'
$wsQ3.Range("G7").Value = '0.2123'
$wsQ3.Range("H7").Value = 10
$wsQ3.Range("A8").Value = 6
$wsQ3.Range("B8").Value = '011845'
$wsQ3.Range("C8").Value = '博时周期优选混合A'
$wsQ3.Range("D8").Value = '2.21'
$wsQ3.Range("E8").Value = '78.30'
$wsQ3.Range("F8").Value = '3.05'
$wsQ3.Range("G8").Value = '0.0674'
$wsQ3.Range("H8").Value = 6
$wsQ3.Range("A9").Value = 7
$wsQ3.Range("B9").Value = '160642'
$wsQ3.Range("C9").Value = '鹏华增瑞灵活配置混合（LOF）'
$wsQ3.Range("D9").Value = '2.05'
$wsQ3.Range("E9").Value = '90.96'
$wsQ3.Range("F9").Value = '3.17'
$wsQ3.Range("G9").Value = '0.0650'
$wsQ3.Range("H9").Value = 9
$wsQ3.Range("A10").Value = 8
$wsQ3.Range("B10").Value = '002556'
$wsQ3.Range("C10").Value = '博时丝路主题股票C'
$wsQ3.Range("D10").Value = '1.35'
$wsQ3.Range("E10").Value = '83.01'
$wsQ3.Range("F10").Value = '2.62'
$wsQ3.Range("G10").Value = '0.0354'
$wsQ3.Range("H10").Value = 7
$wsQ3.Range("A11").Value = 9
$wsQ3.Range("B11").Value = '015031'
$wsQ3.Range("C11").Value = '博时远见回报混合C'
$wsQ3.Range("D11").Value = '1.17'
$wsQ3.Range("E11").Value = '76.26'
$wsQ3.Range("F11").Value = '2.65'
$wsQ3.Range("G11").Value = '0.0310'
$wsQ3.Range("H11").Value = 10
$wsQ3.Range("A12").Value = 10
$wsQ3.Range("B12").Value = '011340'
$wsQ3.Range("C12").Value = '博时战略新材料主题混合A'
$wsQ3.Range("D12").Value = '0.91'
$wsQ3.Range("E12").Value = '79.55'
$wsQ3.Range("F12").Value = '3.30'
$wsQ3.Range("G12").Value = '0.0300'
$wsQ3.Range("H12").Value = 5
$wsQ3.Range("A13").Value = 11
$wsQ3.Range("B13").Value = '014212'
$wsQ3.Range("C13").Value = '博时研究优享混合A'
$wsQ3.Range("D13").Value = '0.80'
$wsQ3.Range("E13").Value = '79.50'
$wsQ3.Range("F13").Value = '2.83'
$wsQ3.Range("G13").Value = '0.0226'
$wsQ3.Range("H13").Value = 8
$wsQ3.Range("A14").Value = 12
$wsQ3.Range("B14").Value = '015030'
$wsQ3.Range("C14").Value = '博时远见回报混合A'
$wsQ3.Range("D14").Value = '0.79'
$wsQ3.Range("E14").Value = '76.26'
$wsQ3.Range("F14").Value = '2.65'
$wsQ3.Range("G14").Value = '0.0209'
$wsQ3.Range("H14").Value = 10
$wsQ3.Range("A15").Value = 13
$wsQ3.Range("B15").Value = '011341'
$wsQ3.Range("C15").Value = '博时战略新材料主题混合C'
$wsQ3.Range("D15").Value = '0.60'
$wsQ3.Range("E15").Value = '79.55'
$wsQ3.Range("F15").Value = '3.30'
$wsQ3.Range("G15").Value = '0.0198'
$wsQ3.Range("H15").Value = 5
$wsQ3.Range("A16").Value = 14
$wsQ3.Range("B16").Value = '009741'
$wsQ3.Range("C16").Value = '博时研究臻选三年持有期灵活配置混合C'
$wsQ3.Range("D16").Value = '0.49'
$wsQ3.Range("E16").Value = '82.29'
$wsQ3.Range("F16").Value = '2.79'
$wsQ3.Range("G16").Value = '0.0137'
$wsQ3.Range("H16").Value = 10
$wsQ3.Range("A17").Value = 15
$wsQ3.Range("B17").Value = '013798'
$wsQ3.Range("C17").Value = '博时优质鑫选一年持有期混合C'
$wsQ3.Range("D17").Value = '0.47'
$wsQ3.Range("E17").Value = '80.56'
$wsQ3.Range("F17").Value = '2.68'
$wsQ3.Range("G17").Value = '0.0126'
$wsQ3.Range("H17").Value = 10
$wsQ3.Range("A18").Value = 16
$wsQ3.Range("B18").Value = '014913'
$wsQ3.Range("C18").Value = '博时研究回报混合A'
$wsQ3.Range("D18").Value = '0.30'
$wsQ3.Range("E18").Value = '78.32'
$wsQ3.Range("F18").Value = '3.48'
$wsQ3.Range("G18").Value = '0.0104'
$wsQ3.Range("H18").Value = 7
$wsQ3.Range("A19").Value = 17
$wsQ3.Range("B19").Value = '014914'
$wsQ3.Range("C19").Value = '博时研究回报混合C'
$wsQ3.Range("D19").Value = '0.10'
$wsQ3.Range("E19").Value = '78.32'
$wsQ3.Range("F19").Value = '3.48'
$wsQ3.Range("G19").Value = '0.0035'
$wsQ3.Range("H19").Value = 7
$wsQ3.Range("A20").Value = 18
$wsQ3.Range("B20").Value = '011846'
$wsQ3.Range("C20").Value = '博时周期优选混合C'
$wsQ3.Range("D20").Value = '0.10'
$wsQ3.Range("E20").Value = '78.30'
$wsQ3.Range("F20").Value = '3.05'
$wsQ3.Range("G20").Value = '0.0030'
$wsQ3.Range("H20").Value = 6
$wsQ3.Range("A21").Value = 19
$wsQ3.Range("B21").Value = '014213'
$wsQ3.Range("C21").Value = '博时研究优享混合C'
$wsQ3.Range("D21").Value = '0.07'
$wsQ3.Range("E21").Value = '79.50'
$wsQ3.Range("F21").Value = '2.83'
$wsQ3.Range("G21").Value = '0.0020'
$wsQ3.Range("H21").Value = 8

# Apply the workbook's existing bold+boxed header style to the header row and to
# the row-index column A, by copying the format from a cell that already has it.
$styleSource = $wsTotal.Range("B1")
$styleSource.Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$styleSource2 = $wsTotal.Range("A2")
$styleSource2.Copy()
$wsQ3.Range("A2:A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2. Add a summary row for "2022-Q3" at the top of the data in "Zong Ji" (总计)
$wsTotal.Rows(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 20
$wsTotal.Range("D2").Value = 12.07

# Re-copy the index-column style onto the new row and renumber the index column
# sequentially (0, 1, 2, ...) now that a row was inserted at the top
$styleSource3 = $wsTotal.Range("A3")
$styleSource3.Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

